# Generate Report for Handback
#
# The bf6e3734-f831-4b2d-bd5f-fb8e8a94aa38.md file has now been handed
# back (it was previously only "Ready for handoff"). Update its status
# to "Handed back: in sync with en-US" on every sheet that reports it,
# and record the new handback timestamps on the per-language sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: one row per localized file, one status column per language ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus   # zh-cn status
$overview.Range("C3").Value = $newStatus   # de-de status

# --- zh-cn sheet: Status (col B) + Latest Handback DateTime (col G) for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("G3").Value = "2016-03-10 03:24:26"

# --- de-de sheet: Status (col B) + Latest Handback DateTime (col G) for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("G3").Value = "2016-03-10 03:24:31"
